# Chapter 3 - Design of Solution: fix up a handful of paragraphs so the
# misspelled / grammar-flagged words sit in their own runs (mirrors what
# Word's proofing pass does when it wraps a flagged span in its own run).
#
# wdFindContinue = 1, wdReplaceNone = 0 (unused here - we only locate text)
$d = $word.ActiveDocument

function Split-Run([string]$phrase) {
    # Re-resolve from the top of the story each time so earlier edits in
    # this script can't shift/stale an old Range object.
    $r = $d.Content
    $ok = $r.Find.Execute($phrase, $true, $false, $false, $false, $false, `
                           $true, 1, $false, $null, 0)
    if (-not $ok) {
        throw "Could not find phrase: $phrase"
    }
    # Flipping a character property on and back off forces Word to break
    # the enclosing run into (before | phrase | after) without altering
    # the phrase's own formatting once the flag is restored to its
    # original (unset) state.
    $r.Font.Bold = $true
    $r.Font.Bold = $false
}

# 3.3 - "...utilizing the Non-Object Oriented technique."
#   -> isolate "Object Oriented" (flagged by the grammar checker)
Split-Run "Object Oriented"

# 3.3 - "...When a program grows in size, it is divided..."
#   -> isolate "grows in size" (flagged by the grammar checker)
Split-Run "grows in size"

# 3.4 - Admin actor bullet: "...manage user permitions in the system."
#   -> isolate "permitions" (flagged by the spell checker)
Split-Run "permitions"

# 3.4 - Delivery Manager actor bullet: "...Manage all deliverys and..."
#   -> isolate "deliverys" (flagged by the spell checker)
Split-Run "deliverys"
